$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This script adds two newly handed-off localization files to the
# "Generate Report for Handoff" localization-status workbook:
#   - 7abf3104-48c9-4ff1-abb4-3fa0aacff4aa.md
#   - c771cdfd-9498-4d79-80a2-d7f330b13b6f.md
# across the Overview, zh-cn and de-de sheets/tables.
# ---------------------------------------------------------------------------

$file1Name    = "7abf3104-48c9-4ff1-abb4-3fa0aacff4aa.md"
$file1Path    = "e2e\7abf3104-48c9-4ff1-abb4-3fa0aacff4aa.md"
$file1ZhXlf   = "7abf3104-48c9-4ff1-abb4-3fa0aacff4aa.25e04ecef63ff027da7fb2ee33130e11e1ef5eae.zh-cn.xlf"
$file1DeXlf   = "7abf3104-48c9-4ff1-abb4-3fa0aacff4aa.25e04ecef63ff027da7fb2ee33130e11e1ef5eae.de-de.xlf"

$file2Name    = "c771cdfd-9498-4d79-80a2-d7f330b13b6f.md"
$file2Path    = "e2e\c771cdfd-9498-4d79-80a2-d7f330b13b6f.md"
$file2ZhXlf   = "c771cdfd-9498-4d79-80a2-d7f330b13b6f.d827e40143dd4aa532e6113c93486bcd6ac53ce4.zh-cn.xlf"
$file2DeXlf   = "c771cdfd-9498-4d79-80a2-d7f330b13b6f.d827e40143dd4aa532e6113c93486bcd6ac53ce4.de-de.xlf"

$status       = "Ready for handoff"
$ext          = ".md"
# A leading apostrophe forces Excel to keep these as literal text instead of
# auto-converting to boolean True/False values.
$trueText     = "'True"
$falseText    = "'False"

$overviewDate = "2017-01-03 04:21:37"
$zhHoDate     = "2017-01-03 04:21:27"
$deHoDate     = "2017-01-03 04:21:37"
$epochDate    = "0001-01-01 00:00:00"

$sourceRepoUrl = "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/363282c018352b6f48372b14ee248381f3a3a76a/e2e/"
$zhRepoUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test1-zhcn/blob/a8c516b6daa377623978e5314824899970b17e04/e2e/"
$deRepoUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test1-dede/blob/12f13263be5988ca993a7771763d56f2299ddbee/e2e/"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

$row = $loOverview.ListRows.Add()
$wsOverview.Range("A3").Value = $file1Name
$wsOverview.Range("B3").Value = $file1Path
$wsOverview.Range("C3").Value = $ext
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = $overviewDate
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $sourceRepoUrl + $file1Name, "", "", $file1Path)

$row = $loOverview.ListRows.Add()
$wsOverview.Range("A4").Value = $file2Name
$wsOverview.Range("B4").Value = $file2Path
$wsOverview.Range("C4").Value = $ext
$wsOverview.Range("E4").Value = $status
$wsOverview.Range("F4").Value = $status
$wsOverview.Range("G4").Value = $overviewDate
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $sourceRepoUrl + $file2Name, "", "", $file2Path)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)

$row = $loZh.ListRows.Add()
$wsZh.Range("A3").Value = $file1Name
$wsZh.Range("B3").Value = $ext
$wsZh.Range("C3").Value = $status
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = $file1ZhXlf
$wsZh.Range("H3").Value = $zhHoDate
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I3").Value = $dupEmpty
$wsZh.Range("J3").Value = $dupEmpty
$wsZh.Range("K3").Value = $dupEmpty
$wsZh.Range("L3").Value = $epochDate
$wsZh.Range("L3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("M3").Value = $dupEmpty
$wsZh.Range("N3").Value = $dupEmpty
$wsZh.Range("O3").Value = "True"
$wsZh.Range("P3").Value = $dupEmpty
$wsZh.Range("Q3").Value = "False"
$wsZh.Range("R3").Value = $dupEmpty
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $zhRepoUrl + $file1Name, "", "", $file1Name)

$row = $loZh.ListRows.Add()
$wsZh.Range("A4").Value = $file2Name
$wsZh.Range("B4").Value = $ext
$wsZh.Range("C4").Value = $status
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "False"
$wsZh.Range("G4").Value = $file2ZhXlf
$wsZh.Range("H4").Value = $zhHoDate
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = $dupEmpty
$wsZh.Range("J4").Value = $dupEmpty
$wsZh.Range("K4").Value = $dupEmpty
$wsZh.Range("L4").Value = $epochDate
$wsZh.Range("L4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("M4").Value = $dupEmpty
$wsZh.Range("N4").Value = $dupEmpty
$wsZh.Range("O4").Value = "True"
$wsZh.Range("P4").Value = $dupEmpty
$wsZh.Range("Q4").Value = "False"
$wsZh.Range("R4").Value = $dupEmpty
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $zhRepoUrl + $file2Name, "", "", $file2Name)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)

$row = $loDe.ListRows.Add()
$wsDe.Range("A3").Value = $file1Name
$wsDe.Range("B3").Value = $ext
$wsDe.Range("C3").Value = $status
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = $file1DeXlf
$wsDe.Range("H3").Value = $deHoDate
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I3").Value = $dupEmpty
$wsDe.Range("J3").Value = $dupEmpty
$wsDe.Range("K3").Value = $dupEmpty
$wsDe.Range("L3").Value = $epochDate
$wsDe.Range("L3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M3").Value = $dupEmpty
$wsDe.Range("N3").Value = $dupEmpty
$wsDe.Range("O3").Value = "True"
$wsDe.Range("P3").Value = $dupEmpty
$wsDe.Range("Q3").Value = "False"
$wsDe.Range("R3").Value = $dupEmpty
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $deRepoUrl + $file1Name, "", "", $file1Name)

$row = $loDe.ListRows.Add()
$wsDe.Range("A4").Value = $file2Name
$wsDe.Range("B4").Value = $ext
$wsDe.Range("C4").Value = $status
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "False"
$wsDe.Range("G4").Value = $file2DeXlf
$wsDe.Range("H4").Value = $deHoDate
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = $dupEmpty
$wsDe.Range("J4").Value = $dupEmpty
$wsDe.Range("K4").Value = $dupEmpty
$wsDe.Range("L4").Value = $epochDate
$wsDe.Range("L4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M4").Value = $dupEmpty
$wsDe.Range("N4").Value = $dupEmpty
$wsDe.Range("O4").Value = "True"
$wsDe.Range("P4").Value = $dupEmpty
$wsDe.Range("Q4").Value = "False"
$wsDe.Range("R4").Value = $dupEmpty
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $deRepoUrl + $file2Name, "", "", $file2Name)
